$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 6: "Volatility, atomicity, and interlocking" -> bold + red
$para6 = $tr.Paragraphs(6, 1)
$para6.Font.Bold = $true
$para6.Font.Color.RGB = 255

# Paragraph 7, run 2: "thread pool and " -> mark dirty (no visible formatting change,
# but bring rPr in line with sibling runs by touching the run's font).
$para7 = $tr.Paragraphs(7, 1)
$runs7 = $para7.Runs()
$runs7.Item(2).Text = $runs7.Item(2).Text
